# Update cryptos list (Thu Mar 14 15:42:49 UTC 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($cellRef, $text)
    # Leading apostrophe forces Excel to keep the value as literal text
    # (preventing auto-conversion of numeric-looking strings to numbers,
    # which would lose meaningful trailing/format zeros). Resetting the
    # style back to "Normal" afterwards avoids leaving a stray
    # text-number-format style on the cell.
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Column D = Price, Column E = Volume(1h)
# Row 2: Bitcoin
Set-PriceText "D2" "70.937.61"
$ws.Range("E2").Value = "  -2.59%  "

# Row 3: Ethereum
Set-PriceText "D3" "3.842.74"
$ws.Range("E3").Value = "  -3.51%  "

# Row 4: TetherUSD (price unchanged)
$ws.Range("E4").Value = "  +0.44%  "

# Row 5: BNB
Set-PriceText "D5" "589.90"
$ws.Range("E5").Value = "  +0.72%  "

# Row 6: Solana
Set-PriceText "D6" "165.81"
$ws.Range("E6").Value = "  +4.68%  "

# Row 7: XRP
Set-PriceText "D7" "0.669"
$ws.Range("E7").Value = "  -1.65%  "

# Row 8: USDC (price unchanged)
$ws.Range("E8").Value = "  +0.23%  "

# Row 9: Cardano
Set-PriceText "D9" "0.748"
$ws.Range("E9").Value = "  -0.39%  "

# Row 10: Dogecoin
Set-PriceText "D10" "0.172"
$ws.Range("E10").Value = "  +2.92%  "

# Row 11: Avalanche
Set-PriceText "D11" "53.13"
$ws.Range("E11").Value = "  -2.14%  "

# Row 12: ShibaInu
Set-PriceText "D12" "0.0000318"
$ws.Range("E12").Value = "  +0.15%  "

# Row 13: Polkadot
Set-PriceText "D13" "11.09"
$ws.Range("E13").Value = "  +2.01%  "

# Row 14: WrappedliquidstakedEther2.0
Set-PriceText "D14" "4.472.67"
$ws.Range("E14").Value = "  -3.02%  "

# Row 15: WrappedEther
Set-PriceText "D15" "3.854.29"
$ws.Range("E15").Value = "  -2.78%  "

# Row 16: Chainlink
Set-PriceText "D16" "20.57"
$ws.Range("E16").Value = "  +0.45%  "

# Row 17: Uniswap
Set-PriceText "D17" "13.74"
$ws.Range("E17").Value = "  -2.10%  "

# Row 18: Polygon (price unchanged)
$ws.Range("E18").Value = "  -6.60%  "

# Row 19: TRON (price unchanged)
$ws.Range("E19").Value = "  -1.97%  "

# Row 20: WrappedBTC
Set-PriceText "D20" "70.914.49"
$ws.Range("E20").Value = "  -2.12%  "

# Row 21: BitcoinCash
Set-PriceText "D21" "430.39"
$ws.Range("E21").Value = "  -0.81%  "

# Row 22: PancakeSwap
Set-PriceText "D22" "4.67"
$ws.Range("E22").Value = "  -0.35%  "

# Row 23: Litecoin
Set-PriceText "D23" "93.73"
$ws.Range("E23").Value = "  -2.41%  "

# Row 24: ImmutableX (price unchanged)
$ws.Range("E24").Value = "  -5.36%  "

# Row 25: InternetComputer(DFINITY)
Set-PriceText "D25" "13.69"
$ws.Range("E25").Value = "  -4.28%  "

# Row 26: Toncoin
Set-PriceText "D26" "4.09"
$ws.Range("E26").Value = "  -7.01%  "

# Row 27: RenderToken
Set-PriceText "D27" "10.78"
$ws.Range("E27").Value = "  -4.65%  "

# Row 28: LEO
Set-PriceText "D28" "5.94"
$ws.Range("E28").Value = "  +0.14%  "

# Row 29: Filecoin
Set-PriceText "D29" "10.10"
$ws.Range("E29").Value = "  -6.49%  "

# Row 30: EthereumClassic
Set-PriceText "D30" "34.74"
$ws.Range("E30").Value = "  -4.61%  "

# Row 31: NEARProtocol
Set-PriceText "D31" "7.79"
$ws.Range("E31").Value = "  -0.47%  "

# Row 32: InjectiveProtocol
Set-PriceText "D32" "49.72"
$ws.Range("E32").Value = "  -2.04%  "

# Row 33: Cosmos
Set-PriceText "D33" "13.42"
$ws.Range("E33").Value = "  -1.53%  "

# Row 34: Hedera (price unchanged)
$ws.Range("E34").Value = "  -5.69%  "

# Row 35: OKB
Set-PriceText "D35" "68.50"
$ws.Range("E35").Value = "  -0.36%  "

# Row 36: PEPE
Set-PriceText "D36" "0.0₃0967"
$ws.Range("E36").Value = "  +12.52%  "

# Row 37: Bittensor
Set-PriceText "D37" "616.74"
$ws.Range("E37").Value = "  -9.09%  "

# Row 38: TheGraph
Set-PriceText "D38" "0.415"
$ws.Range("E38").Value = "  -5.09%  "

# Row 39: Dai (price unchanged)
$ws.Range("E39").Value = "  +0.33%  "

# Row 40: FirstDigitalUSD
Set-PriceText "D40" "0.999"
$ws.Range("E40").Value = "  -0.17%  "

# Row 41: ThetaToken (price unchanged)
$ws.Range("E41").Value = "  -2.15%  "

# Row 42: Kaspa
Set-PriceText "D42" "0.142"
$ws.Range("E42").Value = "  -2.92%  "

# Row 43: dogwifhat
Set-PriceText "D43" "3.20"
$ws.Range("E43").Value = "  +33.55%  "

# Row 44: VeChain (price unchanged)
$ws.Range("E44").Value = "  -4.55%  "

# Row 45: THORChain
Set-PriceText "D45" "10.08"
$ws.Range("E45").Value = "  -7.71%  "

# Rows 46 & 47 swap: Stellar and Fetch.AI trade ranking positions
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-PriceText "D46" "2.62"
$ws.Range("E46").Value = "  -2.69%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-PriceText "D47" "0.143"
$ws.Range("E47").Value = "  -4.10%  "

# Row 48: ApeXProtocol
Set-PriceText "D48" "3.35"
$ws.Range("E48").Value = "  -0.67%  "

# Row 49: Maker
Set-PriceText "D49" "2.828.86"
$ws.Range("E49").Value = "  +2.56%  "

# Row 50: WEMIXToken
Set-PriceText "D50" "2.74"
$ws.Range("E50").Value = "  -18.97%  "

# Row 51: FLOKI
Set-PriceText "D51" "0.000271"
$ws.Range("E51").Value = "  +0.14%  "
